$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: add "*" marker in column A (shared string index 7 -> "*") and move the
# value from C13 to D13 (game became optional/"download-available" like the other
# rows marked with "*").
$ws.Range("A13").Value = "*"
$ws.Range("C13").Value = $null
$ws.Range("D13").Value = 9.58

# Update the active selection to D13, matching the author's last edit location.
$ws.Range("D13").Select()
